# Edit description (per the supplied xml diff / commit message):
#   1. "hello" -> "Hello", split across two runs: "H" and "ello".
#   2. A new paragraph "Thank you for helping me." is added right after
#      the first paragraph. The pre-existing hidden "_GoBack" bookmark
#      (which originally trailed "hello") ends up trailing the new
#      paragraph's run instead, matching the target ordering of
#      <w:r>Thank you for helping me.</w:r> then
#      <w:bookmarkStart/>/<w:bookmarkEnd/>.

$d = $word.ActiveDocument

# Locate "hello" without hard-coding character offsets.
$helloRange = $d.Content
$helloRange.Find.Execute("hello", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$helloStart = $helloRange.Start
$helloEnd = $helloRange.End

# --- Step 1: split the paragraph right after "hello" ------------------
# The bookmark sits at the end of the paragraph, immediately after
# "hello" and before the paragraph mark. Inserting a paragraph mark at
# that exact point pushes the bookmark into a new, second paragraph;
# inserting the new sentence right after the break (but still before
# the now-relocated, zero-width bookmark) keeps the new run in front of
# the bookmark, as in the target.
$splitPoint = $d.Range($helloEnd, $helloEnd)
$splitPoint.InsertBefore("`r")

$newParaStart = $d.Range($helloEnd + 1, $helloEnd + 1)
$newParaStart.InsertBefore("Thank you for helping me.")

# --- Step 2: turn "hello" into two runs, "H" + "ello" (-> "Hello") ----
# A plain Range.Text assignment collapses back into a single run when
# the formatting is identical on both sides, so the split is expressed
# directly as WordprocessingML via InsertXML, which replaces the exact
# range with the literal runs supplied.
$splitRange = $d.Range($helloStart, $helloEnd)
$splitRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>H</w:t></w:r><w:r><w:t>ello</w:t></w:r></w:p>')
